$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Return_accumulate" test rows (31-36) to the test catalogue,
# matching the pattern of the existing rows (Test name / Description / macro).
# Columns A and C are filled first for every new row, then column B is
# filled in afterwards for the last two rows (this mirrors how the
# shared-string table ends up ordered in the authored workbook).

$ws.Cells.Item(31, 1).Value = "Return_accumulate1"
$ws.Cells.Item(31, 2).Value = "Test simple return aggregated monthly"
$ws.Cells.Item(31, 3).Value = "return_accumulate_test1"

$ws.Cells.Item(32, 1).Value = "Return_accumulate2"
$ws.Cells.Item(32, 2).Value = "Test compound return aggregated monthly"
$ws.Cells.Item(32, 3).Value = "return_accumulate_test2"

$ws.Cells.Item(33, 1).Value = "Return_accumulate3"
$ws.Cells.Item(33, 2).Value = "Test simple return aggregated quarterly"
$ws.Cells.Item(33, 3).Value = "return_accumulate_test3"

$ws.Cells.Item(34, 1).Value = "Return_accumulate4"
$ws.Cells.Item(34, 2).Value = "Test compound return aggregated quarterly"
$ws.Cells.Item(34, 3).Value = "return_accumulate_test4"

$ws.Cells.Item(35, 1).Value = "Return_accumulate5"
$ws.Cells.Item(35, 3).Value = "return_accumulate_test5"

$ws.Cells.Item(36, 1).Value = "Return_accumulate6"
$ws.Cells.Item(36, 3).Value = "return_accumulate_test6"

$ws.Cells.Item(35, 2).Value = "Test simple return aggregated yearly"
$ws.Cells.Item(36, 2).Value = "Test compound return aggregated yearly"

# Move the selection to the new last cell, like the author did after
# typing in the final row.
$ws.Range("A36").Select() | Out-Null
